$d = $word.ActiveDocument

# Locate the title text run "How should time be saved in your company?"
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("How should time be saved in your company?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found text, then insert the new " Explain."
    # text there. Because the insertion point sits exactly at the boundary of
    # the existing (bold, Noto Sans, 13.5pt) title run, the new text
    # immediately inherits that same formatting.
    $endPos = $findRange.End
    $insertRange = $d.Range($endPos, $endPos)
    $insertRange.InsertAfter(" Explain.")

    # The newly inserted text now occupies the range right after the original
    # end. Nudge its size away and back to force the engine to keep it as its
    # own distinct run (matching the source diff, which models " Explain." as
    # a second <w:r> with identical run formatting) instead of silently
    # re-absorbing it into the preceding run merely because the formatting
    # already matches.
    $newTextRange = $d.Range($endPos, $endPos + 9)
    $newTextRange.Font.Size = 11
    $newTextRange.Font.NameAscii = "Noto Sans"
    $newTextRange.Font.NameFarEast = "Times New Roman"
    $newTextRange.Font.NameOther = "Noto Sans"
    $newTextRange.Font.NameBi = "Noto Sans"
    $newTextRange.Font.Bold = $true
    $newTextRange.Font.Color = 1511698
    $newTextRange.Font.Size = 13.5
}
